$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the balance banner text
$ws.Range("A9").Value = "BALANCE : 1399.0"

# New transaction rows (14-17)
$ws.Range("A14").Value = "2025-10-21 18:54:23"
$ws.Range("B14").Value = "Deposit"
$ws.Range("C14").Value = 1
$ws.Range("D14").Value = 1301

$ws.Range("A15").Value = "2025-10-21 18:54:59"
$ws.Range("B15").Value = "Withdraw"
$ws.Range("C15").Value = 1
$ws.Range("D15").Value = 1300

$ws.Range("A16").Value = "2025-10-21 19:28:05"
$ws.Range("B16").Value = "Deposit"
$ws.Range("C16").Value = 100
$ws.Range("D16").Value = 1400

$ws.Range("A17").Value = "2025-10-21 19:57:10"
$ws.Range("B17").Value = "Withdraw"
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = 1399
